# Scheduled runner update: refresh market-price-derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 566.8461
$ws.Range("I41").Value = 154.66667
$ws.Range("J41").Value = 920.1429000000001
$ws.Range("K41").Value = 154.66667
$ws.Range("L41").Value = 920.1429000000001
$ws.Range("M41").Value = 285.33333
$ws.Range("N41").Value = -1800.1429
$ws.Range("H76").Value = 4003071.5
$ws.Range("I76").Value = 4350912.5
$ws.Range("K76").Value = 4350912.5
$ws.Range("M76").Value = -4350597.5
$ws.Range("H79").Value = 4003071.5
$ws.Range("I79").Value = 4350912.5
$ws.Range("K79").Value = 4350912.5
$ws.Range("M79").Value = -4349820.5
$ws.Range("H112").Value = 1788
$ws.Range("J112").Value = 1866.9231
$ws.Range("L112").Value = 5600.7693
$ws.Range("N112").Value = -7816.7693
$ws.Range("H132").Value = 1013.52747
$ws.Range("I132").Value = 1035.0122
$ws.Range("J132").Value = 817.7778
$ws.Range("K132").Value = 3105.036599999999
$ws.Range("L132").Value = 2453.3334
$ws.Range("M132").Value = -575.0365999999995
$ws.Range("N132").Value = -7513.3334
$ws.Range("H137").Value = 819.7037
$ws.Range("I137").Value = 672.63635
$ws.Range("J137").Value = 920.8125
$ws.Range("K137").Value = 2017.90905
$ws.Range("L137").Value = 2762.4375
$ws.Range("M137").Value = 532.09095
$ws.Range("N137").Value = -7862.4375
$ws.Range("H138").Value = 1328.04
$ws.Range("I138").Value = 715.17645
$ws.Range("J138").Value = 1965.9183
$ws.Range("K138").Value = 2145.52935
$ws.Range("L138").Value = 5897.7549
$ws.Range("M138").Value = 2994.47065
$ws.Range("N138").Value = -16177.7549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21980.82
$ws.Range("I32").Value = 21692.8
$ws.Range("J32").Value = 23041.947
$ws.Range("K32").Value = 21692.8
$ws.Range("L32").Value = 23041.947
$ws.Range("M32").Value = -21405.8
$ws.Range("N32").Value = -23615.947
$ws.Range("H61").Value = 891.0714
$ws.Range("I61").Value = 804.1892
$ws.Range("J61").Value = 1534
$ws.Range("K61").Value = 804.1892
$ws.Range("L61").Value = 1534
$ws.Range("M61").Value = -592.1892
$ws.Range("N61").Value = -1958
$ws.Range("H74").Value = 910.8958
$ws.Range("I74").Value = 862.73334
$ws.Range("J74").Value = 1633.3334
$ws.Range("K74").Value = 862.73334
$ws.Range("L74").Value = 1633.3334
$ws.Range("M74").Value = 11.26666
$ws.Range("N74").Value = -3381.3334
$ws.Range("H77").Value = 910.8958
$ws.Range("I77").Value = 862.73334
$ws.Range("J77").Value = 1633.3334
$ws.Range("K77").Value = 4313.6667
$ws.Range("L77").Value = 8166.666999999999
$ws.Range("M77").Value = 54.33330000000024
$ws.Range("N77").Value = -16902.667
$ws.Range("H98").Value = 24177.5
$ws.Range("J98").Value = 24177.5
$ws.Range("L98").Value = 24177.5
$ws.Range("N98").Value = -30167.5
$ws.Range("H102").Value = 2153.6428
$ws.Range("I102").Value = 2624
$ws.Range("J102").Value = 977.75
$ws.Range("K102").Value = 2624
$ws.Range("L102").Value = 977.75
$ws.Range("M102").Value = -1002
$ws.Range("N102").Value = -4221.75
$ws.Range("H132").Value = 2034.9678
$ws.Range("I132").Value = 2010.9615
$ws.Range("J132").Value = 2159.8
$ws.Range("K132").Value = 6032.8845
$ws.Range("L132").Value = 6479.400000000001
$ws.Range("M132").Value = -3502.8845
$ws.Range("N132").Value = -11539.4
$ws.Range("H136").Value = 891.0714
$ws.Range("I136").Value = 804.1892
$ws.Range("J136").Value = 1534
$ws.Range("K136").Value = 2412.5676
$ws.Range("L136").Value = 4602
$ws.Range("M136").Value = 137.4323999999997
$ws.Range("N136").Value = -9702

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1275.9286
$ws.Range("I94").Value = 1334.8462
$ws.Range("J94").Value = 510
$ws.Range("K94").Value = 1334.8462
$ws.Range("L94").Value = 510
$ws.Range("M94").Value = -883.8462
$ws.Range("N94").Value = -1412

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1958928
$ws.Range("I31").Value = 2309937.8
$ws.Range("J31").Value = 3302
$ws.Range("K31").Value = 2309937.8
$ws.Range("L31").Value = 3302
$ws.Range("M31").Value = -2309642.8
$ws.Range("N31").Value = -3892
$ws.Range("H34").Value = 1958928
$ws.Range("I34").Value = 2309937.8
$ws.Range("J34").Value = 3302
$ws.Range("K34").Value = 2309937.8
$ws.Range("L34").Value = 3302
$ws.Range("M34").Value = -2309735.8
$ws.Range("N34").Value = -3706
$ws.Range("H52").Value = 13677.25
$ws.Range("I52").Value = 4709
$ws.Range("J52").Value = 16666.666
$ws.Range("K52").Value = 4709
$ws.Range("L52").Value = 16666.666
$ws.Range("M52").Value = -4415
$ws.Range("N52").Value = -17254.666
$ws.Range("H58").Value = 5599.52
$ws.Range("I58").Value = 1928.7693
$ws.Range("J58").Value = 9576.166999999999
$ws.Range("K58").Value = 1928.7693
$ws.Range("L58").Value = 9576.166999999999
$ws.Range("M58").Value = -1725.7693
$ws.Range("N58").Value = -9982.166999999999
$ws.Range("H132").Value = 1671.0428
$ws.Range("I132").Value = 1058.6945
$ws.Range("J132").Value = 2319.4119
$ws.Range("K132").Value = 3176.0835
$ws.Range("L132").Value = 6958.2357
$ws.Range("M132").Value = -646.0835000000002
$ws.Range("N132").Value = -12018.2357
$ws.Range("H134").Value = 1626.683
$ws.Range("I134").Value = 1672.3077
$ws.Range("J134").Value = 1547.6
$ws.Range("K134").Value = 5016.9231
$ws.Range("L134").Value = 4642.799999999999
$ws.Range("M134").Value = -2481.9231
$ws.Range("N134").Value = -9712.799999999999
$ws.Range("H136").Value = 5599.52
$ws.Range("I136").Value = 1928.7693
$ws.Range("J136").Value = 9576.166999999999
$ws.Range("K136").Value = 5786.3079
$ws.Range("L136").Value = 28728.501
$ws.Range("M136").Value = -3236.3079
$ws.Range("N136").Value = -33828.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 505.5965
$ws.Range("I113").Value = 782.875
$ws.Range("J113").Value = 460.32654
$ws.Range("K113").Value = 2348.625
$ws.Range("L113").Value = 1380.97962
$ws.Range("M113").Value = -178.625
$ws.Range("N113").Value = -5720.97962
$ws.Range("H122").Value = 536.65515
$ws.Range("I122").Value = 380.5
$ws.Range("J122").Value = 554.6731
$ws.Range("K122").Value = 3424.5
$ws.Range("L122").Value = 4992.0579
$ws.Range("M122").Value = -974.5
$ws.Range("N122").Value = -9892.0579
$ws.Range("H131").Value = 31777.8
$ws.Range("I131").Value = 112483.445
$ws.Range("J131").Value = 18807.25
$ws.Range("K131").Value = 337450.335
$ws.Range("L131").Value = 56421.75
$ws.Range("M131").Value = -332410.335
$ws.Range("N131").Value = -66501.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 48.076923
$ws.Range("I2").Value = 47.7
$ws.Range("J2").Value = 49.333332
$ws.Range("K2").Value = 47.7
$ws.Range("L2").Value = 49.333332
$ws.Range("M2").Value = 65.3
$ws.Range("N2").Value = -275.333332
$ws.Range("H113").Value = 898.94116
$ws.Range("I113").Value = 305.75
$ws.Range("K113").Value = 305.75
$ws.Range("M113").Value = 1864.25
$ws.Range("H132").Value = 2291.2092
$ws.Range("I132").Value = 2062.3333
$ws.Range("J132").Value = 2677.4375
$ws.Range("K132").Value = 6186.999899999999
$ws.Range("L132").Value = 8032.3125
$ws.Range("M132").Value = -3656.999899999999
$ws.Range("N132").Value = -13092.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1984.491
$ws.Range("I132").Value = 1860.2051
$ws.Range("J132").Value = 2287.4375
$ws.Range("K132").Value = 5580.615299999999
$ws.Range("L132").Value = 6862.3125
$ws.Range("M132").Value = -3050.615299999999
$ws.Range("N132").Value = -11922.3125
$ws.Range("H136").Value = 1895.4231
$ws.Range("I136").Value = 1051.409
$ws.Range("J136").Value = 6537.5
$ws.Range("K136").Value = 3154.227
$ws.Range("L136").Value = 19612.5
$ws.Range("M136").Value = -604.2270000000003
$ws.Range("N136").Value = -24712.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H110").Value = 28000
$ws.Range("J110").Value = 28000
$ws.Range("L110").Value = 28000
$ws.Range("N110").Value = -36180
$ws.Range("H122").Value = 783.6087
$ws.Range("I122").Value = 737
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 2211
$ws.Range("L122").Value = 3015
$ws.Range("M122").Value = 239
$ws.Range("N122").Value = -7915
$ws.Range("H132").Value = 1104.2858
$ws.Range("I132").Value = 945.6667
$ws.Range("K132").Value = 2837.0001
$ws.Range("M132").Value = -307.0001000000002
$ws.Range("H136").Value = 1405.4546
$ws.Range("I136").Value = 1536.5294
$ws.Range("J136").Value = 959.8
$ws.Range("K136").Value = 4609.5882
$ws.Range("L136").Value = 2879.4
$ws.Range("M136").Value = -2059.5882
$ws.Range("N136").Value = -7979.4
